$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '36.814.85'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = "'" + '2.094.12'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'" + '244.87'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').Value = "'" + '0.654'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'" + '54.10'
$ws.Range('E8').Value = '  -5.06%  '
$ws.Range('D9').Value = "'" + '59.08'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').Value = "'" + '0.366'
$ws.Range('E10').Value = '  -4.39%  '
$ws.Range('D11').Value = "'" + '0.0764'
$ws.Range('E11').Value = '  -2.39%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = "'" + '0.928'
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').Value = "'" + '14.94'
$ws.Range('E14').Value = '  -7.85%  '
$ws.Range('D15').Value = "'" + '2.400.01'
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = "'" + '5.48'
$ws.Range('E16').Value = '  -4.50%  '
$ws.Range('D17').Value = "'" + '2.074.24'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = "'" + '36.789.25'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = "'" + '17.12'
$ws.Range('E19').Value = '  -8.98%  '
$ws.Range('D20').Value = "'" + '72.67'
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('D21').Value = "'" + '0.0₃0880'
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').Value = "'" + '5.45'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').Value = "'" + '239.20'
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = "'" + '2.39'
$ws.Range('E25').Value = '  -3.84%  '
$ws.Range('D26').Value = "'" + '9.64'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').Value = "'" + '2.15'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').Value = "'" + '167.01'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('D29').Value = "'" + '21.01'
$ws.Range('E29').Value = '  +4.09%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').Value = "'" + '5.21'
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  +5.04%  '
$ws.Range('D34').Value = "'" + '0.0607'
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('D35').Value = "'" + '2.43'
$ws.Range('E35').Value = '  +7.57%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').Value = "'" + '1.85'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('D38').Value = "'" + '0.0823'
$ws.Range('E38').Value = '  -7.06%  '
$ws.Range('E39').Value = '  -5.38%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = "'" + '4.90'
$ws.Range('E41').Value = '  -7.60%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'" + '0.0220'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').Value = "'" + '0.0960'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('D44').Value = "'" + '96.36'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = "'" + '2.85'
$ws.Range('E45').Value = '  -8.07%  '
$ws.Range('D46').Value = "'" + '7.82'
$ws.Range('E46').Value = '  +14.83%  '
$ws.Range('D47').Value = "'" + '1.415.81'
$ws.Range('E47').Value = '  +11.19%  '
$ws.Range('D48').Value = "'" + '16.03'
$ws.Range('E48').Value = '  -8.84%  '
$ws.Range('D49').Value = "'" + '2.43'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('D50').Value = "'" + '2.90'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('D51').Value = "'" + '2.289.27'
$ws.Range('E51').Value = '  +2.60%  '
